$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct fuzzed/corrupted input values
$ws.Range("E15").Value = 300000000
$ws.Range("E16").Value = -60473972.810000002
$ws.Range("E26").Value = 108613404

# Restore formulas that match the pattern used in the other columns (C,D,F,G)
$ws.Range("E18").Formula = "=SUM(E12:E17)"
$ws.Range("E21").Formula = "=SUM(E18:E20)"

$excel.Calculate()
